{"js": "// The document content was reorganized: several paragraphs/runs swapped their\n// text with each other (a \"rotation\" of text blocks into different slots),\n// while every paragraph's style/formatting stayed exactly where it was.\n//\n// Because several target strings are themselves the *source* string of another\n// replacement (a rotation), we must capture every original value up-front\n// (before any writes) and then apply the new values, so that a later lookup\n// never accidentally matches text that an earlier step just wrote.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Step 1: locate every Range we will need to rewrite, by searching for\n// its current (original) text, and load that text back so we can confirm it\n// before mutating anything. --------------------------------------------------\n\nconst originals = [\n  \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\",\n  \"Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art\",\n  \"11079086 - Herland\u00ed de Souza Andrade\",\n  \"A definir de acordo com o t\u00f3pico programado\",\n  \"To be defined according to the scheduled topic\",\n  \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\",\n  \"Esta disciplina dever\u00e1 conter no m\u00ednimo duas avalia\u00e7\u00f5es denominadas A1 e A2. As avali\u00e7\u00f5es poder\u00e3o ser: escritas, pr\u00e1ticas, semin\u00e1rios, trabalhos de campo, projetos, ou outra forma de avalia\u00e7\u00e3o definida pelo professor.\",\n  \"M\u00e9dia ponderada das avalia\u00e7\u00f5es (M).\",\n  \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 calculada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\",\n  \"Livros, artigos ou texto fornecido pelo docente respons\u00e1vel extra\u00eddos de livros ou revistas especializadas na \u00e1rea de Engenharia de Produ\u00e7\u00e3o.\",\n];\n\nconst ranges = {};\nfor (const text of originals) {\n  const found = body.search(text, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  ranges[text] = found;\n}\nawait context.sync();\n\nfor (const text of originals) {\n  if (ranges[text].items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${text}\" but found ${ranges[text].items.length}`\n    );\n  }\n}\n\n// --- Step 2: the new text that belongs in each of the ranges located above,\n// keyed by the original text so the mapping is unambiguous regardless of\n// write order. -----------------------------------------------------------\n\nconst newTextFor = {\n  \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\":\n    \"A definir de acordo com o t\u00f3pico programado\",\n  \"Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art\":\n    \"To be defined according to the scheduled topic\",\n  \"11079086 - Herland\u00ed de Souza Andrade\":\n    \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\",\n  \"A definir de acordo com o t\u00f3pico programado\":\n    \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\",\n  \"To be defined according to the scheduled topic\":\n    \"Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art\",\n  \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\":\n    \"Esta disciplina dever\u00e1 conter no m\u00ednimo duas avalia\u00e7\u00f5es denominadas A1 e A2. As avali\u00e7\u00f5es poder\u00e3o ser: escritas, pr\u00e1ticas, semin\u00e1rios, trabalhos de campo, projetos, ou outra forma de avalia\u00e7\u00e3o definida pelo professor.\",\n  \"Esta disciplina dever\u00e1 conter no m\u00ednimo duas avalia\u00e7\u00f5es denominadas A1 e A2. As avali\u00e7\u00f5es poder\u00e3o ser: escritas, pr\u00e1ticas, semin\u00e1rios, trabalhos de campo, projetos, ou outra forma de avalia\u00e7\u00e3o definida pelo professor.\":\n    \"M\u00e9dia ponderada das avalia\u00e7\u00f5es (M).\",\n  \"M\u00e9dia ponderada das avalia\u00e7\u00f5es (M).\":\n    \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 calculada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\",\n  \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 calculada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\":\n    \"Livros, artigos ou texto fornecido pelo docente respons\u00e1vel extra\u00eddos de livros ou revistas especializadas na \u00e1rea de Engenharia de Produ\u00e7\u00e3o.\",\n  \"Livros, artigos ou texto fornecido pelo docente respons\u00e1vel extra\u00eddos de livros ou revistas especializadas na \u00e1rea de Engenharia de Produ\u00e7\u00e3o.\":\n    \"11079086 - Herland\u00ed de Souza Andrade\",\n};\n\n// --- Step 3: write the new values into the previously located ranges. Since\n// every Range object was already resolved against the *original* document in\n// step 1, overwriting one range's text cannot affect which range another\n// entry in `originals` points to. -----------------------------------------\n\nfor (const text of originals) {\n  ranges[text].items[0].insertText(newTextFor[text], \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The document content was reorganized: several paragraphs/runs swapped their\n# text with each other (a \"rotation\" of text blocks into different slots),\n# while every paragraph's style/formatting stayed exactly where it was.\n#\n# Because several target strings are themselves the *source* string of another\n# replacement (a rotation), a naive \"find original text, set new text\" loop\n# would, partway through, search for text that an earlier step already wrote\n# over (or, with pre-fetched Range objects, risk writing into a stale\n# position once the document has reflowed). To stay correct and order-\n# independent we do it in two passes with unique placeholder tokens:\n#   Pass 1: find each ORIGINAL string and overwrite it with a unique\n#           placeholder (immediately, one Find+write at a time).\n#   Pass 2: find each placeholder and overwrite it with the real final text.\n# Every Find+write pair below executes back-to-back against a freshly\n# obtained $d.Content range, so no Range object is ever reused after the\n# document has been mutated elsewhere.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Original = \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\"; New = \"A definir de acordo com o t\u00f3pico programado\" },\n    @{ Original = \"Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art\"; New = \"To be defined according to the scheduled topic\" },\n    @{ Original = \"11079086 - Herland\u00ed de Souza Andrade\"; New = \"Complementar a forma\u00e7\u00e3o dos estudantes abordando, com maior profundidade, t\u00f3picos atuais e relevantes e atualizar com temas no estado da arte.\" },\n    @{ Original = \"A definir de acordo com o t\u00f3pico programado\"; New = \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\" },\n    @{ Original = \"To be defined according to the scheduled topic\"; New = \"Complement the training of students by addressing, in greater depth, current and relevant topics and updating them with themes in the state of the art\" },\n    @{ Original = \"O conte\u00fado desta disciplina optativa ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares ao conte\u00fado regular do curso de gradua\u00e7\u00e3o.\"; New = \"Esta disciplina dever\u00e1 conter no m\u00ednimo duas avalia\u00e7\u00f5es denominadas A1 e A2. As avali\u00e7\u00f5es poder\u00e3o ser: escritas, pr\u00e1ticas, semin\u00e1rios, trabalhos de campo, projetos, ou outra forma de avalia\u00e7\u00e3o definida pelo professor.\" },\n    @{ Original = \"Esta disciplina dever\u00e1 conter no m\u00ednimo duas avalia\u00e7\u00f5es denominadas A1 e A2. As avali\u00e7\u00f5es poder\u00e3o ser: escritas, pr\u00e1ticas, semin\u00e1rios, trabalhos de campo, projetos, ou outra forma de avalia\u00e7\u00e3o definida pelo professor.\"; New = \"M\u00e9dia ponderada das avalia\u00e7\u00f5es (M).\" },\n    @{ Original = \"M\u00e9dia ponderada das avalia\u00e7\u00f5es (M).\"; New = \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 calculada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\" },\n    @{ Original = \"A recupera\u00e7\u00e3o ser\u00e1 composta por uma \u00fanica prova (RC) englobando toda a mat\u00e9ria ministrada ao longo do semestre. A m\u00e9dia final, para os alunos em recupera\u00e7\u00e3o, ser\u00e1 calculada com base na rela\u00e7\u00e3o: MF=(M+RC)/2\"; New = \"Livros, artigos ou texto fornecido pelo docente respons\u00e1vel extra\u00eddos de livros ou revistas especializadas na \u00e1rea de Engenharia de Produ\u00e7\u00e3o.\" },\n    @{ Original = \"Livros, artigos ou texto fornecido pelo docente respons\u00e1vel extra\u00eddos de livros ou revistas especializadas na \u00e1rea de Engenharia de Produ\u00e7\u00e3o.\"; New = \"11079086 - Herland\u00ed de Souza Andrade\" }\n)\n\n# Pass 1: original text -> unique placeholder\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $placeholder = \"@@ROTATE_PLACEHOLDER_$i@@\"\n    $rng = $d.Content\n    $found = $rng.Find.Execute($replacements[$i].Original)\n    if (-not $found) {\n        throw \"Could not find expected original text for item $i\"\n    }\n    $rng.Text = $placeholder\n}\n\n# Pass 2: placeholder -> final text\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $placeholder = \"@@ROTATE_PLACEHOLDER_$i@@\"\n    $rng = $d.Content\n    $found = $rng.Find.Execute($placeholder)\n    if (-not $found) {\n        throw \"Could not find placeholder for item $i\"\n    }\n    $rng.Text = $replacements[$i].New\n}\n"}
